$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Emberlee"
$ws.Range("D9").Value = "Do Mocap poses, think of things to code"

$ws.Range("D10").Select()
